# Daily attendance processing - 2026-01-23 00:01:27
# Reorders the "Recorded By" (column G) comma-separated list of recorder
# identities on each data row: whichever of the real user addresses
# (dnasr281@gmail.com, preferred, else backup@backdoor.com) is present in
# the list gets moved to the front, with the remaining entries kept in
# their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Priority order: if an entry in this list is present in the cell's
# comma-separated values, it is promoted to the front of the list.
$priority = @("dnasr281@gmail.com", "backup@backdoor.com")

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value()

    if ($raw -eq $null) {
        continue
    }

    $parts = @($raw -split ",\s*")

    if ($parts.Length -lt 2) {
        continue
    }

    $promote = $null
    foreach ($candidate in $priority) {
        if ($parts -contains $candidate) {
            $promote = $candidate
            break
        }
    }

    if ($promote -eq $null) {
        continue
    }

    $idx = [Array]::IndexOf($parts, $promote)

    $rest = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $idx) {
            $rest = $rest + @($parts[$i])
        }
    }

    $newParts = @($promote) + $rest
    $newValue = $newParts -join ", "

    if ($newValue -ne $raw) {
        $cell.Value = $newValue
    }
}
